$d = $word.ActiveDocument

# Mapping of old font sizes (in points) to new font sizes (in points).
# Derived from the half-point w:sz values in the OOXML diff:
#   32 -> 36 (16pt -> 18pt)  Name header
#   18 -> 20 (9pt  -> 10pt)  Contact info / body text / bullets / dates
#   24 -> 26 (12pt -> 13pt)  Section headers
#   20 -> 22 (10pt -> 11pt)  Overview paragraph
#   22 -> 24 (11pt -> 12pt)  Job titles / degree titles
$sizeMap = @{
    16 = 18
    9  = 10
    12 = 13
    10 = 11
    11 = 12
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    # Exclude the trailing paragraph mark so we don't introduce a
    # <w:rPr><w:sz/></w:rPr> inside <w:pPr> (only the run formatting
    # should change, per the target diff).
    $r.End = $r.End - 1

    $oldSize = $r.Font.Size
    if ($sizeMap.ContainsKey($oldSize)) {
        $r.Font.Size = $sizeMap[$oldSize]
    }
}
